$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Start from a clean sheet (content + formatting + column widths) and
# rebuild the new "libreria" layout: usuarios / libros / autores /
# libro_favorito tables near the top, and the new autores_libros
# relationship table below.
# ---------------------------------------------------------------------------
$ws.Range("A1:Z1").EntireColumn.Delete()

# --- usuarios table : A2:B5 -------------------------------------------------
$ws.Range("A2").Value2 = "usuarios"
$ws.Range("A3").Value2 = "id"
$ws.Range("B3").Value2 = "nombre"
$ws.Range("A4").Value2 = 1
$ws.Range("B4").Value2 = "adrian"
$ws.Range("A5").Value2 = 2
$ws.Range("B5").Value2 = "rodrigo"

# --- libros table : D2:E6 ---------------------------------------------------
$ws.Range("D2").Value2 = "libros"
$ws.Range("D3").Value2 = "id"
$ws.Range("E3").Value2 = "titulo"
$ws.Range("D4").Value2 = 1
$ws.Range("E4").Value2 = "el mago de oz"
$ws.Range("D5").Value2 = 2
$ws.Range("E5").Value2 = "un mundo feliz"
$ws.Range("D6").Value2 = 3
$ws.Range("E6").Value2 = "game of thrones"

# --- autores table : F2:G6 --------------------------------------------------
$ws.Range("F2").Value2 = "autores"
$ws.Range("F3").Value2 = "id"
$ws.Range("G3").Value2 = "nombre"
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = "juanjo de las mercedez"
$ws.Range("F5").Value2 = 2
$ws.Range("G5").Value2 = "maria tres cocos"
$ws.Range("F6").Value2 = 3
$ws.Range("G6").Value2 = "excequiel perez"

# --- libro_favorito table : I2:K6 -------------------------------------------
$ws.Range("I2").Value2 = "libro_favorito"
$ws.Range("I3").Value2 = "id"
$ws.Range("J3").Value2 = "id_user"
$ws.Range("K3").Value2 = "id_libro"
$ws.Range("I4").Value2 = 1
$ws.Range("J4").Value2 = 1
$ws.Range("K4").Value2 = 1
$ws.Range("I5").Value2 = 2
$ws.Range("J5").Value2 = 1
$ws.Range("K5").Value2 = 3
$ws.Range("I6").Value2 = 3
$ws.Range("J6").Value2 = 2
$ws.Range("K6").Value2 = 2

# --- autores_libros table : D9:F13 ------------------------------------------
$ws.Range("D9").Value2 = "autores_libros"
$ws.Range("D10").Value2 = "id"
$ws.Range("E10").Value2 = "id_autor"
$ws.Range("F10").Value2 = "id_libro"
$ws.Range("D11").Value2 = 1
$ws.Range("E11").Value2 = 1
$ws.Range("F11").Value2 = 1
$ws.Range("D12").Value2 = 2
$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("D13").Value2 = 3
$ws.Range("E13").Value2 = 2
$ws.Range("F13").Value2 = 2

# ---------------------------------------------------------------------------
# Header highlight colours (same two-tone scheme as before, just a new
# shade: the "id" highlight moved from yellow to orange, and the foreign
# key highlight moved from red to a darker red).
# ---------------------------------------------------------------------------
$ws.Range("A3").Interior.Color = 0xC0FF
$ws.Range("D3").Interior.Color = 0xC0FF
$ws.Range("F3").Interior.Color = 0xC0FF
$ws.Range("I3").Interior.Color = 0xC0FF
$ws.Range("D10").Interior.Color = 0xC0FF

$ws.Range("J3:K3").Interior.Color = 0xC0
$ws.Range("E10:F10").Interior.Color = 0xC0

# ---------------------------------------------------------------------------
# Column widths: a narrow spacer column C, plus the shifted/resized data
# columns (D..K) that now carry the four tables (values chosen to land as
# close as possible to the original bestFit pixel widths).
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 3.2857142857142856
$ws.Columns.Item(4).ColumnWidth = 12.428571428571429
$ws.Columns.Item(5).ColumnWidth = 19.142857142857142
$ws.Columns.Item(6).ColumnWidth = 6.285714285714286
$ws.Columns.Item(7).ColumnWidth = 19.142857142857142
$ws.Columns.Item(8).ColumnWidth = 11.285714285714286
$ws.Columns.Item(11).ColumnWidth = 11.285714285714286

# ---------------------------------------------------------------------------
# View: zoom to 90% and leave the selection on the new relationship table.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 90
[void]$ws.Range("D10").Select()
